# Apply the commit's edits:
#  - Merge the mixed-formatting runs of each "RECIBO DE PAGAMENTO" paragraph
#    into a single run with a uniform rPr (only sz=24, no bold).
#  - Update company name KFP SERVICE DIGITAL LTDA -> GRAVATAI SERVICE DIGITAL LTDA
#  - Update CNPJ 41.230.154/0001-57 -> 00.111.222/0001-33
#  - Update the date line to "Cachoeirinha, 2025-04-04 00:00:00."

$d = $word.ActiveDocument
$brChar = [char]11

# Work from the end backwards so earlier paragraph Start/End offsets
# remain valid while we rewrite text further down the document.
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $rng = $p.Range
    $txt = $rng.Text

    if ($txt.StartsWith("RECIBO DE PAGAMENTO")) {
        # NOTE: the text returned by Range.Text represents each <w:br/> as
        # Chr(11); the original run boundary right after "RECIBO DE
        # PAGAMENTO" is itself a <w:br/>, so it must be consumed explicitly
        # instead of being swallowed into the (lazy) name-capture group.
        $pattern = '^RECIBO DE PAGAMENTO' + $brChar + '(.*?), inscrito\(a\) no CPF sob o n. (\d+), declaro para os devidos fins ter recebido nesta data, da empresa .*?, inscrita no CNPJ sob o n. [\d\.\-/]+, a import.ncia de R\$([\d\.]+) concernente ao pagamento de um domingo trabalhado\.'
        if ($txt -match $pattern) {
            $name = $matches[1]
            $cpf = $matches[2]
            $value = $matches[3]

            $newText = "RECIBO DE PAGAMENTO" + $brChar + `
                $name + ", inscrito(a) no CPF sob o nº " + $cpf + `
                ", declaro para os devidos fins ter recebido nesta data, da empresa GRAVATAI SERVICE DIGITAL LTDA, inscrita no CNPJ sob o nº 00.111.222/0001-33, a importância de R$" + $value + `
                " concernente ao pagamento de um domingo trabalhado." + $brChar + $brChar + `
                "Cachoeirinha, 2025-04-04 00:00:00." + $brChar + $brChar + `
                "_________________________________________________" + $brChar + `
                "Assinatura" + $brChar

            # Exclude the trailing paragraph mark from the replaced range.
            $paraRange = $d.Range($rng.Start, $rng.End - 1)
            $paraRange.Text = $newText

            $newRange = $d.Range($rng.Start, $rng.Start + $newText.Length)
            $newRange.Font.Size = 12
        }
    }
}
